# Generate Report for Handoff
#
# The "cca6a516-04d2-4af1-a9c8-24c9febfa8cc.md" file (last data row, row 7,
# on every sheet) has just been handed off for localization. Refresh the
# handoff timestamps that the report shows for it:
#
#   Overview!G7 (Latest HO Xliff Generate Date) -> 2016-08-21 22:51:35
#   zh-cn!H7    (Latest Handoff Datetime)        -> 2016-08-21 22:51:30
#   de-de!H7    (Latest Handoff Datetime)        -> 2016-08-21 22:51:35

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-21 22:51:35"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-21 22:51:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-21 22:51:35"
